$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.906.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.169.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.31%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.46%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.168.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.30%  "

$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("E10").Value = "  +6.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.49%  "

$ws.Range("E13").Value = "  +17.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.688.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.985.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.182.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.88%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.32%  "

$ws.Range("E19").Value = "  +1.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.16%  "

$ws.Range("E21").Value = "  +6.00%  "

$ws.Range("E22").Value = "  +7.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.60%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.43%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  +5.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.26%  "

$ws.Range("E35").Value = "  +6.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0897"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "472.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0420"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("E41").Value = "  +5.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.066.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.39%  "

$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("E44").Value = "  +6.49%  "

$ws.Range("E45").Value = "  +8.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0605"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +19.01%  "

$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("E50").Value = "  +8.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.56%  "
